# bowler_data_Smat_form_score.xlsx
# Re-sort the existing bowler rows (2-7) alphabetically by player name and
# append two freshly-scraped bowlers (XC Bartlett, E Malinga) as rows 8-9.
# Also drops the stale AutoFilter sort-state and moves the saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param(
        [string]$range,
        [object[]]$values
    )
    $arr = New-Object 'object[,]' 1,$values.Length
    for ($i = 0; $i -lt $values.Length; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $ws.Range($range).Value = $arr
}

# Rows 3 (DS Rathi), 4 (Gurjapneet Singh) and 7 (Zeeshan Ansari) keep the same
# data/position in the re-sorted table, so they are intentionally left alone.

# Row 2 -> Ashwani Kumar (previously row 6)
Set-Row "A2:Q2" @("Ashwani Kumar", 2, 2, 7, 0, 52, 1, 52, 7.42, 42, 0, 0, "BOWL", 9, 1, 0, -4.810821197000001)

# Row 5 -> Prince Yadav (previously row 2)
Set-Row "A5:Q5" @("Prince Yadav", 8, 8, 27.1, 0, 205, 11, 18.63, 7.54, 14.8, 0, 0, "BOWL", 10, 1, 0, 0.1179517521000002)

# Row 6 -> PVSN Raju (previously row 5)
Set-Row "A6:Q6" @("PVSN Raju", 7, 7, 22.5, 0, 188, 7, 26.85, 8.23, 19.5, 0, 0, "BOWL", 11, 1, 0, -1.51207914)

# New rows 8-9 need the same direct formatting (shaded fill) as the rest of
# the data rows, so stamp that in first by copying an existing row's format.
$ws.Range("A7:Q7").Copy()
$ws.Range("A8:Q9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 8 -> XC Bartlett (new)
Set-Row "A8:Q8" @("XC Bartlett", 9, 9, 35.1, 0, 333, 12, 27.75, 9.4600000000000009, 17.5, 1, 0, "BOWL", 8, 1, 1, -0.60655534550000001)

# Row 9 -> E Malinga (new)
Set-Row "A9:Q9" @("E Malinga", 1, 1, 1.5, 0, 26, 0, 100, 14.18, 100, 0, 0, "BOWL", 9, 1, 0, -9.7912563190000004)

# The data is no longer maintained under an explicit Consistency_Score sort,
# so drop the worksheet's remembered sort state.
$ws.Sort.SortFields.Clear()

# Leave the selection where the author last left it.
$null = $ws.Range("L16").Select()
